# Requirements plastic vanger.docx - add two new requirement rows
# ("De micro controller moet de sensoren kunnen uitlezen." and
#  "De micro controller heeft een wifi chip.") to the requirements table,
# right before the "De motor moet altijd op de minimale snelheid draaien."
# row.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Locate the row that should follow the two new rows by matching on its
# first-cell text (robust against row-index drift).
$targetRowIndex = -1
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $cellText = $tbl.Cell($i, 1).Range.Text
    if ($cellText -like "*De motor moet altijd op de*snelheid draaien*") {
        $targetRowIndex = $i
        break
    }
}

if ($targetRowIndex -eq -1) {
    throw "Could not locate anchor row 'De motor moet altijd op de ... snelheid draaien.'"
}

$anchorRow = $tbl.Rows.Item($targetRowIndex)

# Insert the "wifi chip" row directly above the anchor row first ...
$wifiRow = $tbl.Rows.Add($anchorRow)
$wifiRow.Cells.Item(1).Range.Text = "De micro controller heeft een wifi chip."
$wifiRow.Cells.Item(2).Range.Text = "Ik wil dat de micro controller een wifi chip heeft."
$wifiRow.Cells.Item(4).Range.Text = "De micro controller moet een wifi chip bevatten."
$wifiRow.Cells.Item(5).Range.Text = "Could "

# ... then insert the "sensoren uitlezen" row above that, so the final
# order is: sensoren-row, wifi-row, anchor-row.
$sensorRow = $tbl.Rows.Add($wifiRow)
$sensorRow.Cells.Item(1).Range.Text = "De micro controller moet de sensoren kunnen uitlezen."
$sensorRow.Cells.Item(2).Range.Text = "Ik wil dat de micro controller de sensoren uit kunnen lezen en daarmee de motor kunnen aansturen."
$sensorRow.Cells.Item(4).Range.Text = "De micro controller kan de sensoren uitlezen."
$sensorRow.Cells.Item(5).Range.Text = "Must"
